# Kayıt silindi: 11349401
# The record with Kayıt No 11349401 (sheet1 row 1293 / "Merkez İlçe" row 754)
# is removed; all subsequent rows shift up by one on both the master
# "Kayitlar" sheet and the per-birim "Merkez İlçe" sheet.

$wb = $excel.ActiveWorkbook

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows.Item(1293).Delete()

$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")
$wsMerkez.Rows.Item(754).Delete()
